$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19: fill in values for A19:D19
$ws.Range("A19").Value = "2.2"
$ws.Range("B19").Value = "Add log features"
$ws.Range("C19").Value = "yes"
$ws.Range("D19").Value = 0.505

# Row 20: add new row
$ws.Range("A20").Value = 2.3
$ws.Range("B20").Value = "Change stratified kfold to repeated k fold and increase train set"
$ws.Range("C20").Value = "yes"
$ws.Range("D20").Value = 0.509

$ws.Range("E18").Select()
